$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05003266666666667
$ws.Range("H2").Value = 0.150098
$ws.Range("O2").Value = 0.1810521476743106
$ws.Range("P2").Value = 0.1810521476743105
$ws.Range("Q2").Value = 0.0006174197842222222
$ws.Range("R2").Value = 0.005556778058
$ws.Range("S2").Value = 0.1810521476743106
$ws.Range("T2").Value = 0.1810521476743105

# Add new row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Wnt7b"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05003266666666667
$ws.Range("H3").Value = 0.150098
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05581866666666666
$ws.Range("N3").Value = 0.167456
$ws.Range("O3").Value = 0.8189478523256895
$ws.Range("P3").Value = 0.8189478523256895
$ws.Range("Q3").Value = 0.002792756743111111
$ws.Range("R3").Value = 0.025134810688
$ws.Range("S3").Value = 0.8189478523256895
$ws.Range("T3").Value = 0.8189478523256895
